$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "44.089.31"
Set-TextValue "E2" "  +0.48%  "
Set-TextValue "D3" "2.266.21"
Set-TextValue "E3" "  -0.75%  "
Set-TextValue "E4" "  +0.14%  "
Set-TextValue "D5" "233.12"
Set-TextValue "E5" "  +0.56%  "
Set-TextValue "D6" "0.651"
Set-TextValue "E6" "  +3.36%  "
Set-TextValue "D7" "63.60"
Set-TextValue "E7" "  -1.16%  "
Set-TextValue "E8" "  +0.01%  "
Set-TextValue "D9" "0.451"
Set-TextValue "E9" "  +5.93%  "
Set-TextValue "D10" "0.0978"
Set-TextValue "E10" "  -0.55%  "
Set-TextValue "D11" "58.35"
Set-TextValue "E11" "  +1.28%  "
Set-TextValue "D12" "26.79"
Set-TextValue "E12" "  +1.17%  "
Set-TextValue "E13" "  +2.04%  "
Set-TextValue "D14" "2.602.45"
Set-TextValue "E14" "  -0.77%  "
Set-TextValue "D15" "15.72"
Set-TextValue "E15" "  -0.76%  "
Set-TextValue "D16" "6.17"
Set-TextValue "E16" "  +3.76%  "
Set-TextValue "D17" "0.840"
Set-TextValue "E17" "  +2.23%  "
Set-TextValue "D18" "2.259.80"
Set-TextValue "E18" "  -1.70%  "
Set-TextValue "D19" "43.930.49"
Set-TextValue "E19" "  +0.54%  "
Set-TextValue "D20" "0.0₃0987"
Set-TextValue "E20" "  +2.94%  "
Set-TextValue "D21" "74.01"
Set-TextValue "E21" "  +0.65%  "
Set-TextValue "E22" "  -0.79%  "
Set-TextValue "D23" "248.08"
Set-TextValue "E23" "  -1.06%  "
Set-TextValue "E24" "  +0.02%  "
Set-TextValue "E25" "  -3.87%  "
Set-TextValue "D26" "2.30"
Set-TextValue "E26" "  -1.44%  "
Set-TextValue "D27" "3.33"
Set-TextValue "E27" "  +19.49%  "
Set-TextValue "B28" "Cosmos"
Set-TextValue "C28" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D28" "9.96"
Set-TextValue "E28" "  -0.09%  "
Set-TextValue "B29" "EthereumClassic"
Set-TextValue "C29" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D29" "22.32"
Set-TextValue "E29" "  +8.25%  "
Set-TextValue "D30" "173.57"
Set-TextValue "E30" "  +0.76%  "
Set-TextValue "E31" "  -0.19%  "
Set-TextValue "D32" "1.46"
Set-TextValue "E32" "  +0.19%  "
Set-TextValue "E33" "  +2.09%  "
Set-TextValue "D34" "5.03"
Set-TextValue "E34" "  +6.06%  "
Set-TextValue "D35" "0.0685"
Set-TextValue "E35" "  -1.75%  "
Set-TextValue "D36" "4.99"
Set-TextValue "E36" "  -3.90%  "
Set-TextValue "D37" "3.69"
Set-TextValue "E37" "  -3.16%  "
Set-TextValue "D38" "6.48"
Set-TextValue "E38" "  -5.29%  "
Set-TextValue "D39" "2.31"
Set-TextValue "E39" "  -2.02%  "
Set-TextValue "E40" "  +2.74%  "
Set-TextValue "E41" "  +0.18%  "
Set-TextValue "D42" "8.76"
Set-TextValue "E42" "  +3.31%  "
Set-TextValue "D43" "0.000221"
Set-TextValue "E43" "  -2.21%  "
Set-TextValue "D44" "17.41"
Set-TextValue "E44" "  +2.67%  "
Set-TextValue "E45" "  +0.42%  "
Set-TextValue "E46" "  -2.00%  "
Set-TextValue "B47" "FTXToken"
Set-TextValue "C47" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D47" "4.42"
Set-TextValue "E47" "  -4.06%  "
Set-TextValue "B48" "Cronos"
Set-TextValue "C48" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D48" "0.0951"
Set-TextValue "E48" "  -2.16%  "
Set-TextValue "D49" "2.36"
Set-TextValue "E49" "  +1.63%  "
Set-TextValue "D50" "1.451.58"
Set-TextValue "E50" "  -2.58%  "
Set-TextValue "D51" "9.99"
Set-TextValue "E51" "  -9.52%  "
